$wb = $excel.ActiveWorkbook

# Row 6 on sheet ALC (hunk 0: -926,22 +926,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3667
$ws.Range("I6").Value = 3500.5
$ws.Range("K6").Value = 10501.5
$ws.Range("M6").Value = -10389.5

# Row 17 on sheet ALC (hunk 1: -1465,20 +1465,23)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2833.3333
$ws.Range("I17").Value = 2500
$ws.Range("K17").Value = 7500
$ws.Range("M17").Value = -7332

# Row 33 on sheet ALC (hunk 2: -2237,22 +2240,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 291.17648
$ws.Range("I33").Value = 291.17648
$ws.Range("K33").Value = 291.17648
$ws.Range("M33").Value = -62.17648000000003

# Row 38 on sheet ALC (hunk 3: -2476,22 +2479,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 759
$ws.Range("I38").Value = 97.333336
$ws.Range("K38").Value = 292.000008
$ws.Range("M38").Value = 79.99999200000002

# Row 39 on sheet ALC (hunk 4: -2528,25 +2531,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 101.75
$ws.Range("I39").Value = 37.333332
$ws.Range("J39").Value = 295
$ws.Range("K39").Value = 111.999996
$ws.Range("L39").Value = 885
$ws.Range("M39").Value = 184.000004
$ws.Range("N39").Value = -1477

# Row 58 on sheet ALC (hunk 5: -3450,25 +3453,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1892.7142
$ws.Range("I58").Value = 49.5
$ws.Range("J58").Value = 2630
$ws.Range("K58").Value = 148.5
$ws.Range("L58").Value = 7890
$ws.Range("M58").Value = 1.5
$ws.Range("N58").Value = -8190

# Row 70 on sheet ALC (hunk 6: -4011,25 +4014,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3474.75
$ws.Range("J70").Value = 3466.3333
$ws.Range("L70").Value = 10398.9999
$ws.Range("N70").Value = -10938.9999

# Row 73 on sheet ALC (hunk 7: -4158,25 +4161,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3474.75
$ws.Range("J73").Value = 3466.3333
$ws.Range("L73").Value = 10398.9999
$ws.Range("N73").Value = -12270.9999

# Row 80 on sheet ALC (hunk 8: -4492,25 +4495,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1306.25
$ws.Range("I80").Value = 1416.6666
$ws.Range("J80").Value = 1240
$ws.Range("K80").Value = 4249.9998
$ws.Range("L80").Value = 3720
$ws.Range("M80").Value = -3251.9998
$ws.Range("N80").Value = -5716

# Row 82 on sheet ALC (hunk 9: -4590,22 +4593,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 999.6667
$ws.Range("I82").Value = 999.6667
$ws.Range("K82").Value = 2999.0001
$ws.Range("M82").Value = -2593.0001

# Row 83 on sheet ALC (hunk 10: -4639,25 +4642,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1306.25
$ws.Range("I83").Value = 1416.6666
$ws.Range("J83").Value = 1240
$ws.Range("K83").Value = 12749.9994
$ws.Range("L83").Value = 11160
$ws.Range("M83").Value = -7757.999400000001
$ws.Range("N83").Value = -21144

# Row 85 on sheet ALC (hunk 11: -4737,22 +4740,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 999.6667
$ws.Range("I85").Value = 999.6667
$ws.Range("K85").Value = 2999.0001
$ws.Range("M85").Value = -1595.0001

# Row 99 on sheet ALC (hunk 12: -5438,22 +5441,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 300
$ws.Range("I99").Value = 200
$ws.Range("K99").Value = 600
$ws.Range("M99").Value = 898

# Row 101 on sheet ALC (hunk 13: -5542,22 +5545,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 493
$ws.Range("I101").Value = 366.25
$ws.Range("K101").Value = 1098.75
$ws.Range("M101").Value = 523.25

# Row 112 on sheet ALC (hunk 14: -6078,22 +6081,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2490.196
$ws.Range("J112").Value = 2490.196
$ws.Range("L112").Value = 7470.588
$ws.Range("N112").Value = -9686.588

# Row 118 on sheet ALC (hunk 15: -6363,22 +6366,19)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

# Row 127 on sheet ALC (hunk 16: -6792,22 +6792,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 498.5
$ws.Range("I127").Value = 498.5
$ws.Range("K127").Value = 1495.5
$ws.Range("M127").Value = 3464.5

# Row 45 on sheet ARM (hunk 17: -9695,22 +9695,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1452.8113
$ws.Range("I45").Value = 1411.7451
$ws.Range("K45").Value = 1411.7451
$ws.Range("M45").Value = -1034.7451

# Row 61 on sheet ARM (hunk 18: -10470,19 +10470,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1500
$ws.Range("J61").Value = 1500
$ws.Range("L61").Value = 1500
$ws.Range("N61").Value = -1924

# Row 110 on sheet ARM (hunk 19: -12796,22 +12799,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 656.2941
$ws.Range("I110").Value = 642.44446
$ws.Range("K110").Value = 642.44446
$ws.Range("M110").Value = 1402.55554

# Row 122 on sheet ARM (hunk 20: -13372,22 +13375,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6239.5
$ws.Range("I122").Value = 7102.4
$ws.Range("K122").Value = 21307.2
$ws.Range("M122").Value = -18857.2

# Row 136 on sheet ARM (hunk 21: -14034,19 +14037,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1500
$ws.Range("J136").Value = 1500
$ws.Range("L136").Value = 4500
$ws.Range("N136").Value = -9600

# Row 107 on sheet BSM (hunk 22: -19453,25 +19459,25)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4468.6875
$ws.Range("I107").Value = 4466.6665
$ws.Range("J107").Value = 4499
$ws.Range("K107").Value = 4466.6665
$ws.Range("L107").Value = 4499
$ws.Range("M107").Value = -2546.6665
$ws.Range("N107").Value = -8339

# Row 134 on sheet BSM (hunk 23: -20734,19 +20740,22)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1000
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = -465

# Row 10 on sheet CRP (hunk 24: -21609,25 +21618,25)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 3968.6667
$ws.Range("I10").Value = 951.75
$ws.Range("J10").Value = 10002.5
$ws.Range("K10").Value = 951.75
$ws.Range("L10").Value = 10002.5
$ws.Range("M10").Value = -812.75
$ws.Range("N10").Value = -10280.5

# Row 59 on sheet CRP (hunk 25: -24004,22 +24013,22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 53999
$ws.Range("J59").Value = 53999
$ws.Range("L59").Value = 53999
$ws.Range("N59").Value = -56289

# Row 107 on sheet CRP (hunk 26: -26317,19 +26326,22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 996.3333
$ws.Range("I107").Value = 994.5
$ws.Range("K107").Value = 994.5
$ws.Range("M107").Value = 925.5

# Row 2 on sheet CUL (hunk 27: -28090,22 +28102,22)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 38.285713
$ws.Range("I2").Value = 24.666666
$ws.Range("K2").Value = 147.999996
$ws.Range("M2").Value = -34.99999600000001

# Row 86 on sheet CUL (hunk 28: -32227,22 +32239,22)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1901
$ws.Range("I86").Value = 602
$ws.Range("K86").Value = 1806
$ws.Range("M86").Value = -620

# Row 89 on sheet CUL (hunk 29: -32371,22 +32383,22)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1901
$ws.Range("I89").Value = 602
$ws.Range("K89").Value = 5418
$ws.Range("M89").Value = 510

# Row 129 on sheet CUL (hunk 30: -34292,25 +34304,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 962.8
$ws.Range("J129").Value = 1046
$ws.Range("L129").Value = 3138
$ws.Range("N129").Value = -13138

# Row 131 on sheet CUL (hunk 31: -34393,25 +34405,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2387
$ws.Range("I131").Value = 2422
$ws.Range("J131").Value = 2343.25
$ws.Range("K131").Value = 7266
$ws.Range("L131").Value = 7029.75
$ws.Range("M131").Value = -2226
$ws.Range("N131").Value = -17109.75

# Row 63 on sheet GSM (hunk 32: -37940,19 +37952,22)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 75000
$ws.Range("J63").Value = 75000
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76372

# Row 66 on sheet GSM (hunk 33: -38078,19 +38093,22)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 75000
$ws.Range("J66").Value = 75000
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -231864

# Row 70 on sheet GSM (hunk 34: -38265,22 +38283,22)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1472
$ws.Range("I70").Value = 1472
$ws.Range("K70").Value = 1472
$ws.Range("M70").Value = -1202

# Row 73 on sheet GSM (hunk 35: -38409,22 +38427,22)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 1472
$ws.Range("I73").Value = 1472
$ws.Range("K73").Value = 1472
$ws.Range("M73").Value = -536

# Row 102 on sheet GSM (hunk 36: -39791,22 +39809,22)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 200749.5
$ws.Range("I102").Value = 200749.5
$ws.Range("K102").Value = 200749.5
$ws.Range("M102").Value = -199127.5

# Row 26 on sheet LTW (hunk 37: -42922,25 +42940,25)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 4564.3335
$ws.Range("I26").Value = 6004.5
$ws.Range("J26").Value = 1684
$ws.Range("K26").Value = 6004.5
$ws.Range("L26").Value = 1684
$ws.Range("M26").Value = -5709.5
$ws.Range("N26").Value = -2274

# Row 61 on sheet LTW (hunk 38: -44598,22 +44616,22)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 10000
$ws.Range("I61").Value = 10000
$ws.Range("K61").Value = 10000
$ws.Range("M61").Value = -9798

# Row 113 on sheet LTW (hunk 39: -47071,22 +47089,22)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 10000
$ws.Range("K113").Value = 10000
$ws.Range("M113").Value = -7830

# Row 62 on sheet WVR (hunk 40: -51367,25 +51385,22)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4248

# Row 65 on sheet WVR (hunk 41: -51517,25 +51532,22)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -21240

# Row 132 on sheet WVR (hunk 42: -54737,22 +54749,22)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1727
$ws.Range("I132").Value = 1727
$ws.Range("K132").Value = 5181
$ws.Range("M132").Value = -2651

# Row 138 on sheet WVR (hunk 43: -55028,19 +55040,22)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 120000
$ws.Range("J138").Value = 120000
$ws.Range("L138").Value = 120000
$ws.Range("N138").Value = -130280
